$d = $word.ActiveDocument

# 1) Battery bank description: "Daewoo Deep cycle" -> "[bn]"
$d.Content.Find.Execute("Daewoo Deep cycle", $true, $false, $false, $false, $false, $true, 1, $false, "[bn]", 2)

# 2) Battery bank spec: "180 AH& 12 Vdc" -> "[bs]"
$d.Content.Find.Execute("180 AH& 12 Vdc", $true, $false, $false, $false, $false, $true, 1, $false, "[bs]", 2)

# 3) Payment clause: merge "Greevo (Pvt.) Limited." run pair
$d.Content.Find.Execute("Greevo (Pvt.) Limited.", $true, $false, $false, $false, $false, $true, 1, $false, "Greevo (Pvt.) Limited.", 2)

# 4) Payment clause: merge trailing "Greevo (Pvt." run pair
$d.Content.Find.Execute("by Greevo (Pvt.", $true, $false, $false, $false, $false, $true, 1, $false, "by Greevo (Pvt.", 2)

# 5) Delayed payments clause
$d.Content.Find.Execute("M/S Greevo Pvt Ltd reserves", $true, $false, $false, $false, $false, $true, 1, $false, "M/S Greevo Pvt Ltd reserves", 2)

# 6) Modifications clause
$d.Content.Find.Execute("consent from Greevo will nullify", $true, $false, $false, $false, $false, $true, 1, $false, "consent from Greevo will nullify", 2)

# 7) Net metering clause (two Greevo mentions in one paragraph)
$d.Content.Find.Execute("Greevo will assist Customer in procuring Net Metering. The obligation to procure required approvals rests with the Customer. In no circumstances shall Greevo be held liable on any account", $true, $false, $false, $false, $false, $true, 1, $false, "Greevo will assist Customer in procuring Net Metering. The obligation to procure required approvals rests with the Customer. In no circumstances shall Greevo be held liable on any account", 2)

# 8) NEPRA clause
$d.Content.Find.Execute("NEPRA. Greevo will assist Customer to prepare", $true, $false, $false, $false, $false, $true, 1, $false, "NEPRA. Greevo will assist Customer to prepare", 2)

# 9) Signature block: "Greevo PVT Ltd"
$d.Content.Find.Execute("Greevo PVT Ltd", $true, $false, $false, $false, $false, $true, 1, $false, "Greevo PVT Ltd", 2)
